# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 2023-11-03 (serial 45233) to 2023-11-13 (serial 45243).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
